$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last refreshed" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 12:42"

# --- Austria (row 16) ---
$ws.Range("B16").Value = 7269
$ws.Range("C16").Value = 360
$ws.Range("E16").Value = 6986

# --- Noruega (row 20) ---
$ws.Range("B20").Value = 3443
$ws.Range("C20").Value = 71
$ws.Range("E20").Value = 3421
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 16

# --- Eslovenia (row 49) ---
$ws.Range("E49").Value = 613
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 9

# --- Barein (row 60) ---
$ws.Range("B60").Value = 466
$ws.Range("C60").Value = 8
$ws.Range("E60").Value = 252

# --- Irak moves ahead of Serbia: rows 61-63 become Irak / Serbia / Libano ---
# Row 61: Serbia -> Irak (new data)
$ws.Range("A61").Value = "Irak"
$ws.Range("B61").Value = 458
$ws.Range("C61").Value = 76
$ws.Range("D61").Value = 122
$ws.Range("E61").Value = 296
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 40

# Row 62: Libano -> Serbia (carries the old Serbia data, unchanged)
$ws.Range("A62").Value = "Serbia"
$ws.Range("B62").Value = 457
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 15
$ws.Range("E62").Value = 435
$ws.Range("F62").Value = 21
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 7

# Row 63: Irak -> Libano (carries the old Libano data, unchanged)
$ws.Range("A63").Value = "Libano"
$ws.Range("B63").Value = 391
$ws.Range("C63").Value = 23
$ws.Range("D63").Value = 23
$ws.Range("E63").Value = 361
$ws.Range("F63").Value = 3
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 7

# --- Albania (row 85) ---
$ws.Range("D85").Value = 31
$ws.Range("E85").Value = 147

# --- Moldavia (row 86) ---
$ws.Range("F86").Value = 33

# --- Vietnam (row 87) ---
$ws.Range("B87").Value = 163
$ws.Range("C87").Value = 10
$ws.Range("E87").Value = 143

# --- Malta moves ahead of Ghana: rows 91-93 become Malta / Ghana / Reunion ---
# Row 91: Ghana -> Malta (new data)
$ws.Range("A91").Value = "Malta"
$ws.Range("B91").Value = 139
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 2
$ws.Range("E91").Value = 137
$ws.Range("H91").Value = 0

# Row 92: Reunion -> Ghana (carries old Ghana data, unchanged)
$ws.Range("A92").Value = "Ghana"
$ws.Range("B92").Value = 136
$ws.Range("C92").Value = 4
$ws.Range("E92").Value = 131
$ws.Range("F92").Value = 1
$ws.Range("H92").Value = 4

# Row 93: Malta -> Reunion (carries old Reunion data, unchanged)
$ws.Range("A93").Value = "Reunion"
$ws.Range("B93").Value = 135
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = 134
$ws.Range("F93").Value = 0

# --- Camboya (row 101) ---
$ws.Range("B101").Value = 99
$ws.Range("C101").Value = 1
$ws.Range("E101").Value = 88

# --- Bielorrusia (row 104) ---
$ws.Range("D104").Value = 32
$ws.Range("E104").Value = 62

# --- Madagascar moves ahead of Barbados: rows 138-139 become Madagascar / Barbados ---
# Row 138: Barbados -> Madagascar (new data)
$ws.Range("A138").Value = "Madagascar"
$ws.Range("C138").Value = 1

# Row 139: Madagascar -> Barbados (carries old Barbados data, unchanged)
$ws.Range("A139").Value = "Barbados"
$ws.Range("B139").Value = 24
$ws.Range("E139").Value = 24

# --- Maldivas moves ahead of El Salvador: rows 146-148 become Maldivas / El Salvador / Tanzania ---
# Row 146: El Salvador -> Maldivas (new data)
$ws.Range("A146").Value = "Maldivas"
$ws.Range("B146").Value = 14
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 9
$ws.Range("E146").Value = 5

# Row 147: Tanzania -> El Salvador (carries old El Salvador data, unchanged)
$ws.Range("A147").Value = "El Salvador"

# Row 148: Maldivas -> Tanzania (carries old Tanzania data, unchanged)
$ws.Range("A148").Value = "Tanzania"
$ws.Range("D148").Value = 0
$ws.Range("E148").Value = 13
